# Auto-generated edit script: updates cached profit-calculation values
# across ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1005.9167
$ws.Range("I40").Value = 985.8333
$ws.Range("K40").Value = 985.8333
$ws.Range("M40").Value = -810.8333

# Row 105
$ws.Range("H105").Value = 64335
$ws.Range("J105").Value = 64335
$ws.Range("L105").Value = 64335
$ws.Range("N105").Value = -71323

# Row 121
$ws.Range("H121").Value = 1127.4
$ws.Range("J121").Value = 1495.6666
$ws.Range("L121").Value = 4486.9998
$ws.Range("N121").Value = -7980.9998

# Row 137
$ws.Range("H137").Value = 938.7143
$ws.Range("I137").Value = 614.46155
$ws.Range("K137").Value = 1843.38465
$ws.Range("M137").Value = 706.61535

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1111999.8
$ws.Range("J2").Value = 1500
$ws.Range("L2").Value = 1500
$ws.Range("N2").Value = -1726

# Row 32
$ws.Range("H32").Value = 3831.3142
$ws.Range("I32").Value = 3148.082
$ws.Range("J32").Value = 8462.111000000001
$ws.Range("K32").Value = 3148.082
$ws.Range("L32").Value = 8462.111000000001
$ws.Range("M32").Value = -2861.082
$ws.Range("N32").Value = -9036.111000000001

# Row 61
$ws.Range("H61").Value = 1308.1389
$ws.Range("I61").Value = 724.4231
$ws.Range("K61").Value = 724.4231
$ws.Range("M61").Value = -512.4231

# Row 74
$ws.Range("H74").Value = 1166.9333
$ws.Range("I74").Value = 878.9143
$ws.Range("J74").Value = 2175
$ws.Range("K74").Value = 878.9143
$ws.Range("L74").Value = 2175
$ws.Range("M74").Value = -4.914300000000026
$ws.Range("N74").Value = -3923

# Row 77
$ws.Range("H77").Value = 1166.9333
$ws.Range("I77").Value = 878.9143
$ws.Range("J77").Value = 2175
$ws.Range("K77").Value = 4394.5715
$ws.Range("L77").Value = 10875
$ws.Range("M77").Value = -26.57150000000001
$ws.Range("N77").Value = -19611

# Row 110
$ws.Range("H110").Value = 1271.697
$ws.Range("I110").Value = 1001.9286
$ws.Range("K110").Value = 1001.9286
$ws.Range("M110").Value = 1043.0714

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").ClearContents() | Out-Null

# Row 116
$ws.Range("H116").Value = 1111999.8
$ws.Range("J116").Value = 1500
$ws.Range("L116").Value = 1500
$ws.Range("N116").Value = -6088

# Row 132
$ws.Range("H132").Value = 1229.3585
$ws.Range("J132").Value = 1843.75
$ws.Range("L132").Value = 5531.25
$ws.Range("N132").Value = -10591.25

# Row 136
$ws.Range("H136").Value = 1308.1389
$ws.Range("I136").Value = 724.4231
$ws.Range("K136").Value = 2173.2693
$ws.Range("M136").Value = 376.7307000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1111999.8
$ws.Range("J3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("N3").Value = -1728

# Row 20
$ws.Range("H20").Value = 2167.5
$ws.Range("I20").Value = 1958.9231
$ws.Range("K20").Value = 1958.9231
$ws.Range("M20").Value = -1711.9231

# Row 134
$ws.Range("H134").Value = 9822.682000000001
$ws.Range("I134").Value = 9616.611000000001
$ws.Range("J134").Value = 10750
$ws.Range("K134").Value = 28849.833
$ws.Range("L134").Value = 32250
$ws.Range("M134").Value = -26314.833
$ws.Range("N134").Value = -37320

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1677.7778
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50

# Row 31
$ws.Range("H31").Value = 1460.3
$ws.Range("I31").Value = 1058.6207
$ws.Range("K31").Value = 1058.6207
$ws.Range("M31").Value = -763.6206999999999

# Row 34
$ws.Range("H34").Value = 1460.3
$ws.Range("I34").Value = 1058.6207
$ws.Range("K34").Value = 1058.6207
$ws.Range("M34").Value = -856.6206999999999

# Row 58
$ws.Range("H58").Value = 4686.875
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents() | Out-Null

# Row 95
$ws.Range("H95").Value = 26033.2
$ws.Range("J95").Value = 26033.2
$ws.Range("L95").Value = 26033.2
$ws.Range("N95").Value = -31525.2

# Row 122
$ws.Range("H122").Value = 4714.8184
$ws.Range("I122").Value = 3470.6667
$ws.Range("J122").Value = 6207.8
$ws.Range("K122").Value = 10412.0001
$ws.Range("L122").Value = 18623.4
$ws.Range("M122").Value = -7962.000100000001
$ws.Range("N122").Value = -23523.4

# Row 132
$ws.Range("H132").Value = 1620.9395
$ws.Range("I132").Value = 1028.7826
$ws.Range("K132").Value = 3086.3478
$ws.Range("M132").Value = -556.3478

# Row 136
$ws.Range("H136").Value = 4686.875
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 73.333336
$ws.Range("I8").Value = 73.333336
$ws.Range("K8").Value = 220.000008
$ws.Range("M8").Value = -81.00000800000001

# Row 26
$ws.Range("H26").Value = 1042.2
$ws.Range("J26").Value = 900
$ws.Range("L26").Value = 2700
$ws.Range("N26").Value = -3276

# Row 38
$ws.Range("H38").Value = 474
$ws.Range("I38").Value = 197.5
$ws.Range("K38").Value = 592.5
$ws.Range("M38").Value = -245.5

# Row 59
$ws.Range("H59").Value = 300
$ws.Range("J59").Value = 100
$ws.Range("L59").Value = 300
$ws.Range("N59").Value = -1380

# Row 107
$ws.Range("H107").Value = 420.8
$ws.Range("J107").Value = 489.7143
$ws.Range("L107").Value = 1469.1429
$ws.Range("N107").Value = -5309.1429

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 350.2143
$ws.Range("I55").Value = 301.72223
$ws.Range("J55").Value = 437.5
$ws.Range("K55").Value = 301.72223
$ws.Range("L55").Value = 437.5
$ws.Range("M55").Value = -128.72223
$ws.Range("N55").Value = -783.5

# Row 104
$ws.Range("H104").Value = 12775.5
$ws.Range("J104").Value = 12775.5
$ws.Range("L104").Value = 12775.5
$ws.Range("N104").Value = -19763.5

# Row 136
$ws.Range("H136").Value = 2102.5476
$ws.Range("I136").Value = 1340.4193
$ws.Range("K136").Value = 4021.2579
$ws.Range("M136").Value = -1471.2579

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 15783.857
$ws.Range("J69").Value = 15783.857
$ws.Range("L69").Value = 15783.857
$ws.Range("N69").Value = -17281.857

# Row 72
$ws.Range("H72").Value = 15783.857
$ws.Range("J72").Value = 15783.857
$ws.Range("L72").Value = 47351.571
$ws.Range("N72").Value = -54839.571

# Row 95
$ws.Range("H95").Value = 99992.5
$ws.Range("J95").Value = 99992.5
$ws.Range("L95").Value = 99992.5
$ws.Range("N95").Value = -105484.5

# Row 105
$ws.Range("H105").Value = 44953
$ws.Range("J105").Value = 44953
$ws.Range("L105").Value = 44953
$ws.Range("N105").Value = -51941

# Row 136
$ws.Range("H136").Value = 15433826
$ws.Range("I136").Value = 17362616
$ws.Range("K136").Value = 52087848
$ws.Range("M136").Value = -52085298
